# Apply edits described by the diff.
$wb = $excel.ActiveWorkbook

$wsJaana = $wb.Worksheets.Item("Jaana")
$wsJarno = $wb.Worksheets.Item("Jarno")

# --- Jaana sheet: fill in rows 7 and 8 with new data ---
# Reuse the existing date cell format already used on the "Jarno" sheet
# for the same kind of cell (A7/A8), so no new style gets created.
$wsJarno.Range("A7").Copy()
$wsJaana.Range("A7").PasteSpecial(-4122)  # xlPasteFormats
$wsJaana.Range("A7").Value = 44970
$wsJaana.Range("B7").Value = 0.5
$wsJaana.Range("C7").Value = "Tuntiseurantapohja, backlogin pohtimista"

$wsJarno.Range("A8").Copy()
$wsJaana.Range("A8").PasteSpecial(-4122)  # xlPasteFormats
$wsJaana.Range("A8").Value = 44974
$wsJaana.Range("B8").Value = 0.75
$wsJaana.Range("C8").Value = "Planning, dokumentaatio"

# --- Selection / active sheet updates ---
$wsJaana.Activate()
$wsJaana.Range("C8").Select()

$wb.Save()
